$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-parsed as numbers by Excel (e.g. '240.82'), so they stay text
# just like the original inline-string cells.
$textForceRows = @(
5, 6, 8, 9, 10, 11, 13, 14, 16, 18, 20, 21, 22, 23, 25, 26, 27, 28, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 41, 42, 43, 45, 46, 47, 48, 49, 50, 51
)
foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Update Coin / Link / Price / Volume(1h) columns for each changed row
# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '29.428.32'
# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '1.852.17'
$ws.Cells.Item(3, 5).Value = '  -0.01%  '
# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  +0.21%  '
# Row 5: BNB
$ws.Cells.Item(5, 4).Value = '240.82'
$ws.Cells.Item(5, 5).Value = '  -0.02%  '
# Row 6: XRP
$ws.Cells.Item(6, 4).Value = '0.6297'
$ws.Cells.Item(6, 5).Value = '  -0.54%  '
# Row 7: USDC
$ws.Cells.Item(7, 5).Value = '  +0.12%  '
# Row 8: Dogecoin
$ws.Cells.Item(8, 4).Value = '0.07671'
$ws.Cells.Item(8, 5).Value = '  +1.29%  '
# Row 9: Cardano
$ws.Cells.Item(9, 4).Value = '0.2940'
$ws.Cells.Item(9, 5).Value = '  -0.79%  '
# Row 10: Solana
$ws.Cells.Item(10, 4).Value = '24.56'
$ws.Cells.Item(10, 5).Value = '  -0.40%  '
# Row 11: TRON
$ws.Cells.Item(11, 4).Value = '0.07751'
$ws.Cells.Item(11, 5).Value = '  +0.53%  '
# Row 12: WrappedEther
$ws.Cells.Item(12, 4).Value = '1.854.54'
$ws.Cells.Item(12, 5).Value = '  +0.33%  '
# Row 13: ShibaInu
$ws.Cells.Item(13, 4).Value = '0.00001096'
$ws.Cells.Item(13, 5).Value = '  +8.83%  '
# Row 14: Polkadot
$ws.Cells.Item(14, 4).Value = '5.029'
$ws.Cells.Item(14, 5).Value = '  +0.49%  '
# Row 15: Polygon
$ws.Cells.Item(15, 5).Value = '  -0.89%  '
# Row 16: Litecoin
$ws.Cells.Item(16, 4).Value = '83.60'
$ws.Cells.Item(16, 5).Value = '  +0.22%  '
# Row 17: WrappedliquidstakedEther2.0
$ws.Cells.Item(17, 4).Value = '2.106.26'
$ws.Cells.Item(17, 5).Value = '  +0.18%  '
# Row 18: Uniswap
$ws.Cells.Item(18, 4).Value = '6.149'
$ws.Cells.Item(18, 5).Value = '  -0.15%  '
# Row 19: WrappedBTC
$ws.Cells.Item(19, 4).Value = '29.456.36'
$ws.Cells.Item(19, 5).Value = '  -0.13%  '
# Row 20: BitcoinCash
$ws.Cells.Item(20, 4).Value = '229.29'
$ws.Cells.Item(20, 5).Value = '  +0.06%  '
# Row 21: Avalanche
$ws.Cells.Item(21, 4).Value = '12.46'
$ws.Cells.Item(21, 5).Value = '  -0.43%  '
# Row 22: Dai
$ws.Cells.Item(22, 4).Value = '1.001'
$ws.Cells.Item(22, 5).Value = '  +0.10%  '
# Row 23: Chainlink
$ws.Cells.Item(23, 4).Value = '7.452'
$ws.Cells.Item(23, 5).Value = '  -1.17%  '
# Row 24: BinanceUSD
$ws.Cells.Item(24, 5).Value = '  +0.11%  '
# Row 25: Monero
$ws.Cells.Item(25, 4).Value = '156.74'
$ws.Cells.Item(25, 5).Value = '  -0.10%  '
# Row 26: Stellar
$ws.Cells.Item(26, 4).Value = '0.1387'
$ws.Cells.Item(26, 5).Value = '  -0.96%  '
# Row 27: Cosmos
$ws.Cells.Item(27, 4).Value = '8.393'
$ws.Cells.Item(27, 5).Value = '  -0.07%  '
# Row 28: EthereumClassic
$ws.Cells.Item(28, 4).Value = '17.69'
$ws.Cells.Item(28, 5).Value = '  -0.06%  '
# Row 29: Toncoin
$ws.Cells.Item(29, 5).Value = '  +3.42%  '
# Row 30: PancakeSwap
$ws.Cells.Item(30, 4).Value = '1.468'
$ws.Cells.Item(30, 5).Value = '  +0.08%  '
# Row 31: Hedera
$ws.Cells.Item(31, 4).Value = '0.05724'
$ws.Cells.Item(31, 5).Value = '  +0.51%  '
# Row 32: Filecoin
$ws.Cells.Item(32, 4).Value = '4.133'
$ws.Cells.Item(32, 5).Value = '  -0.02%  '
# Row 33: InternetComputer(DFINITY)
$ws.Cells.Item(33, 4).Value = '4.049'
$ws.Cells.Item(33, 5).Value = '  +0.25%  '
# Row 34: LidoDAOToken
$ws.Cells.Item(34, 4).Value = '1.850'
$ws.Cells.Item(34, 5).Value = '  +0.10%  '
# Row 35: ARBITRUM
$ws.Cells.Item(35, 4).Value = '1.162'
$ws.Cells.Item(35, 5).Value = '  +0.20%  '
# Row 36: ImmutableX
$ws.Cells.Item(36, 4).Value = '0.7053'
$ws.Cells.Item(36, 5).Value = '  -1.48%  '
# Row 37: HuobiToken
$ws.Cells.Item(37, 4).Value = '2.586'
$ws.Cells.Item(37, 5).Value = '  -0.08%  '
# Row 38: MXToken
$ws.Cells.Item(38, 4).Value = '2.783'
$ws.Cells.Item(38, 5).Value = '  +0.18%  '
# Row 39: VeChain
$ws.Cells.Item(39, 4).Value = '0.01795'
$ws.Cells.Item(39, 5).Value = '  -0.77%  '
# Row 40: Maker
$ws.Cells.Item(40, 4).Value = '1.219.41'
$ws.Cells.Item(40, 5).Value = '  -2.35%  '
# Row 41: FraxShare
$ws.Cells.Item(41, 4).Value = '6.507'
$ws.Cells.Item(41, 5).Value = '  +4.78%  '
# Row 42: TrustWalletToken
$ws.Cells.Item(42, 4).Value = '0.9101'
$ws.Cells.Item(42, 5).Value = '  +0.06%  '
# Row 43: PaxDollar
$ws.Cells.Item(43, 4).Value = '1.001'
$ws.Cells.Item(43, 5).Value = '  +0.08%  '
# Row 44: RocketPoolETH
$ws.Cells.Item(44, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(44, 4).Value = '2.014.97'
$ws.Cells.Item(44, 5).Value = '  +0.17%  '
# Row 45: Quant
$ws.Cells.Item(45, 2).Value = 'Quant'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(45, 4).Value = '101.73'
$ws.Cells.Item(45, 5).Value = '  +0.06%  '
# Row 46: Aave
$ws.Cells.Item(46, 2).Value = 'Aave'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(46, 4).Value = '66.40'
$ws.Cells.Item(46, 5).Value = '  +0.35%  '
# Row 47: BabyDogeCoin
$ws.Cells.Item(47, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(47, 4).Value = '0.00000000119'
$ws.Cells.Item(47, 5).Value = '  -0.39%  '
# Row 48: Aptos
$ws.Cells.Item(48, 2).Value = 'Aptos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(48, 4).Value = '7.134'
$ws.Cells.Item(48, 5).Value = '  +0.55%  '
# Row 49: TheSandbox
$ws.Cells.Item(49, 2).Value = 'TheSandbox'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(49, 4).Value = '0.4018'
$ws.Cells.Item(49, 5).Value = '  -0.41%  '
# Row 50: EnergySwap
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '9.005'
$ws.Cells.Item(50, 5).Value = '  -1.09%  '
# Row 51: RenderToken
$ws.Cells.Item(51, 2).Value = 'RenderToken'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(51, 4).Value = '1.682'
$ws.Cells.Item(51, 5).Value = '  -0.62%  '
